$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 1.768647483540576
    "D2" = 0.01306468109396519
    "E2" = 1.254454410282079
    "F2" = 0.4207376222180699
    "G2" = 0.002368925570076562
    "L2" = 0.8549488623162915
    "O2" = 1.303034954123888
    "B3" = 1.636809724457521
    "D3" = 0.01159173079111042
    "E3" = 1.144640531064397
    "F3" = 0.4171397892705926
    "G3" = 0.0023725547546041
    "L3" = 0.7511826128100836
    "O3" = 1.306268060040253
    "B4" = 1.555848833617347
    "D4" = 0.01068607563622237
    "E4" = 1.077432815501766
    "F4" = 0.4155398450059948
    "G4" = 0.002374898240683515
    "L4" = 0.6872320586698493
    "O4" = 1.310336985678077
    "B5" = 1.52285433251069
    "D5" = 0.01031672921270399
    "E5" = 1.050103026458686
    "F5" = 0.4150398023411697
    "G5" = 0.002375882278242808
    "L5" = 0.6611128533875501
    "O5" = 1.312514751756567
    "B6" = 1.517375506009103
    "D6" = 0.0102553832843526
    "E6" = 1.04556852204621
    "F6" = 0.414965911956692
    "G6" = 0.002376047433941953
    "L6" = 0.6567722543272794
    "O6" = 1.312907626230697
    "B7" = 1.555403866469476
    "D7" = 0.01068109561159503
    "E7" = 1.077063997964302
    "F7" = 0.4155324876546302
    "G7" = 0.002374911393969643
    "L7" = 0.6868800428690065
    "O7" = 1.310364257536037
    "B8" = 1.723192684674018
    "D8" = 0.01255708706644754
    "E8" = 1.216546961223202
    "F8" = 0.4193699236983335
    "G8" = 0.002370153068713556
    "L8" = 0.8192200897525765
    "O8" = 1.303714628382977
    "B9" = 2.052119155879382
    "D9" = 0.01622476230484438
    "E9" = 1.49169411178363
    "F9" = 0.4317845066952799
    "G9" = 0.002361731456869596
    "L9" = 1.076827810570478
    "O9" = 1.307405281101126
    "B10" = 2.293730119485531
    "D10" = 0.01891139713546863
    "E10" = 1.694710289645229
    "F10" = 0.443965139946485
    "G10" = 0.00235609255876025
    "L10" = 1.26491687341786
    "O10" = 1.32058979470284
    "B11" = 2.403640098074163
    "D11" = 0.02013166033739822
    "E11" = 1.787232778584467
    "F11" = 0.4501878171948022
    "G11" = 0.002353645093650047
    "L11" = 1.350228553547993
    "O11" = 1.328920728301114
    "B12" = 2.445260175988835
    "D12" = 0.02059344640868943
    "E12" = 1.822290776007833
    "F12" = 0.4526435486367717
    "G12" = 0.002352735129906258
    "L12" = 1.382497466833286
    "O12" = 1.332415675667391
    "B13" = 2.436296573495952
    "D13" = 0.02049400633900689
    "E13" = 1.814739482831527
    "F13" = 0.4521102238261108
    "G13" = 0.002352930359409176
    "L13" = 1.375549424358326
    "O13" = 1.33164777092324
    "B14" = 2.407064218097617
    "D14" = 0.02016965794577885
    "E14" = 1.790116600112668
    "F14" = 0.450387853252991
    "G14" = 0.00235356989328536
    "L14" = 1.352884077000965
    "O14" = 1.329201414686764
    "B15" = 2.389158486639872
    "D15" = 0.01997094516900688
    "E15" = 1.775037128421985
    "F15" = 0.449345827493687
    "G15" = 0.002353963816599802
    "L15" = 1.338996089993771
    "O15" = 1.327747394530974
    "B16" = 2.286547203603959
    "D16" = 0.01883160941687834
    "E16" = 1.688666909914076
    "F16" = 0.443572302560014
    "G16" = 0.00235625486718944
    "L16" = 1.259336444645896
    "O16" = 1.320092699918007
    "B17" = 2.223598154645003
    "D17" = 0.01813215689846714
    "E17" = 1.635723070696116
    "F17" = 0.4402059564513223
    "G17" = 0.002357690436171368
    "L17" = 1.210403046836291
    "O17" = 1.315997724927541
    "B18" = 2.187391542216858
    "D18" = 0.01772967301787531
    "E18" = 1.605287270771925
    "F18" = 0.4383338042023297
    "G18" = 0.002358527220967663
    "L18" = 1.182234263611747
    "O18" = 1.313861738390926
    "B19" = 2.175132622594845
    "D19" = 0.01759336930976474
    "E19" = 1.594985077961923
    "F19" = 0.4377108913924843
    "G19" = 0.002358812448110925
    "L19" = 1.172692768310753
    "O19" = 1.313176058146695
    "B20" = 2.230299190171877
    "D20" = 0.01820663336736317
    "E20" = 1.641357386185149
    "F20" = 0.4405576680592844
    "G20" = 0.002357536470981137
    "L20" = 1.215614536831765
    "O20" = 1.316410907780551
    "B21" = 2.415650484680782
    "D21" = 0.0202649353196378
    "E21" = 1.797348368044908
    "F21" = 0.4508910487123075
    "G21" = 0.002353381590488926
    "L21" = 1.359542441890142
    "O21" = 1.32991069842555
    "B22" = 2.536785523764649
    "D22" = 0.02160838739740711
    "E22" = 1.899423057474934
    "F22" = 0.4582241435132346
    "G22" = 0.002350764242845495
    "L22" = 1.453393066562228
    "O22" = 1.340718734831711
    "B23" = 2.472133849767943
    "D23" = 0.02089153233118424
    "E23" = 1.844933225701794
    "F23" = 0.4542568607296857
    "G23" = 0.002352152222200418
    "L23" = 1.403323062572554
    "O23" = 1.3347670927688
    "B24" = 2.227269702944284
    "D24" = 0.01817296367310206
    "E24" = 1.638810104736137
    "F24" = 0.4403984624660211
    "G24" = 0.002357606042964623
    "L24" = 1.213258535301406
    "O24" = 1.316223428285554
    "B25" = 1.963146132920542
    "D25" = 0.01523388166376094
    "E25" = 1.417099027822076
    "F25" = 0.4278944255521679
    "G25" = 0.002363912980844457
    "L25" = 1.007343999713839
    "O25" = 1.304587447384591
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}
